$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All edited cells hold text data (coin names, links, formatted
# price/percentage strings). Force text number-format first so
# Excel does not auto-convert numeric-looking strings (e.g. "1.001")
# into real numbers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '22.465.18'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.03%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.567.70'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.00%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'

$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.11%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '288.18'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.39%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3704'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.47%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '48.13'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -3.66%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3304'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.14%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07580'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.21%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.133'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.27%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.001'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.10%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.61'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.98%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.919'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.11%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.877'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.18%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.565.93'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.16%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001122'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.13%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06782'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.65%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '87.76'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.64%  '

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.10%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.346'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.47%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '16.50'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +2.08%  '

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.28%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '22.425.42'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.07%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.388'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.13%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.576'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.95%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '154.60'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +3.16%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.73'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.26%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.014'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.17%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '124.49'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.35%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.741.01'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.28%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.069'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.63%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.005'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.03%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.106'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.87%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.776'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.97%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08368'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.85%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02471'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.48%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2250'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.17%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06418'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.41%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.286'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -3.67%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.343'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.48%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.31'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.14%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6264'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.92%  '

$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.92'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.18%  '

$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6119'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +5.86%  '

$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.773'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.02%  '

$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.061'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.48%  '

$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '126.02'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.09%  '

$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'EOS'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.207'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.28%  '

$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07214'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.42%  '

$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '76.70'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.80%  '
